$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new columns before the existing A:B, shifting the current
# nameAr/nameEn columns from A:B to C:D (keeps their data + column-C
# "bestFit" width untouched since it simply rides along with the shift).
$ws.Range("A:B").Insert()

# New "department"/"section" header row and "hr"/"employment" data row
# in the freshly inserted A:B columns (only rows 1-2 get values).
$ws.Range("A1").Value = "department"
$ws.Range("B1").Value = "section"
$ws.Range("A2").Value = "hr"
$ws.Range("B2").Value = "employment"

# Column widths for the two new columns (column C/old-A keeps its original
# bestFit width automatically because it was never touched).
$ws.Columns.Item(1).ColumnWidth = 18.17
$ws.Columns.Item(2).ColumnWidth = 13.83

# Match the saved selection/active cell.
$ws.Range("D1").Select()
